# Update the "想去人数" (column F) counts for specific exhibition/event rows
# on both the "展览" sheet and the "全部类型" sheet, matching the regenerated
# gh-pages data snapshot (commit 456a3b4).
#
# Row numbers differ slightly between the two sheets because "全部类型"
# contains one extra row (a "演出" entry) that is not present on "展览".
# Values are matched/verified against the previous (old) value before
# writing the new one, so the correct row is always targeted even if
# row positions were to shift.

$wb = $excel.ActiveWorkbook

# oldValue -> newValue for column F, applied per-sheet below.
$sheetUpdates = @{
    "展览"     = @{
        2  = @{ Old = 7601; New = 7610 }
        3  = @{ Old = 70;   New = 71 }
        4  = @{ Old = 217;  New = 218 }
        5  = @{ Old = 20;   New = 23 }
        6  = @{ Old = 264;  New = 268 }
        7  = @{ Old = 1143; New = 1144 }
        9  = @{ Old = 20;   New = 21 }
        10 = @{ Old = 144;  New = 148 }
    }
    "全部类型" = @{
        2  = @{ Old = 7601; New = 7610 }
        3  = @{ Old = 70;   New = 71 }
        4  = @{ Old = 217;  New = 218 }
        5  = @{ Old = 20;   New = 23 }
        6  = @{ Old = 264;  New = 268 }
        7  = @{ Old = 1143; New = 1144 }
        10 = @{ Old = 20;   New = 21 }
        11 = @{ Old = 144;  New = 148 }
    }
}

foreach ($sheetName in $sheetUpdates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rowMap = $sheetUpdates[$sheetName]

    foreach ($r in $rowMap.Keys) {
        $cell = $ws.Cells.Item($r, 6)
        $expectedOld = $rowMap[$r].Old
        $newValue = $rowMap[$r].New

        $current = $cell.Value2
        if ($current -eq $expectedOld) {
            $cell.Value = $newValue
        } else {
            Write-Output "WARNING: $sheetName!F$r expected $expectedOld but found $current"
        }
    }
}
